# Fix the mislabeled "Chile Cheese Squares" recipe (and its related
# ingredient lists that used "Chile" instead of "Chili"), and populate
# the previously-empty serving_pax / cuisine / rec_type columns for the
# first two recipes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Repopulate serving_pax / cuisine / rec_type for rows 2 and 3 ---
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = "Western"
$ws.Range("F2").Value = "Dinner"

$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "Western"
$ws.Range("F3").Value = "Snack"

# --- Fix the buggy "Chile Cheese Squares" recipe (row 4): "Chile" -> "Chili" ---
$ws.Range("B4").Value = "Preheat oven to 350 degrees F (175 degrees C).**In a 9x13 inch baking dish, place 1/2 of the shredded Cheddar cheese. Top with green chili peppers. Top the  chilis with the remaining cheese.**In a blender, mix eggs with the juice from the chili peppers. Pour the egg mixture over the cheese mixture.**Bake in the preheated oven 30 to 40 minutes. Cool before cutting into squares.**"
$ws.Range("A4").Value = "Chili Cheese Squares  "
$ws.Range("G4").Value = "Cheddar,Chili Pepper,Egg"

# --- Fix the related ingredient lists elsewhere that also said "Chile" ---
$ws.Range("G24").Value = "Potato,Bell Pepper,Chili,Salt,Black Pepper,Paprika,Vegetable Oil,Water"
$ws.Range("G56").Value = "Pea,Pepper,Onion,Chili Pepper,Pimento,Garlic,Salad Dressing"

# --- Move the active selection to B3 ---
$ws.Range("B3").Select()

# --- Page setup (paper size / orientation) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
